# --- Step 1: rename Sheet1 -> Articles, insert ID column ---
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Articles"

$ws1.Columns.Item(1).Insert()
$ws1.Range("A1").Value = "ID"
$ws1.Range("A2").Value = 1
$ws1.Range("A3").Value = 2

# --- Step 2b: create Topics sheet with headers ID, Name ---
$wsT = $wb.Worksheets.Add($null, $ws1)
$wsT.Name = "Topics"

# --- Step 2: create Articles_Topics sheet with headers ID, Article ID, Topic ID, Description ---
$wsAT = $wb.Worksheets.Add($null, $wsT)
$wsAT.Name = "Articles_Topics"
$wsAT.Range("A1").Value = "ID"
$wsAT.Range("B1").Value = "Article ID"
$wsAT.Range("C1").Value = "Topic ID"
$wsAT.Range("D1").Value = "Description"

$wsT.Range("A1").Value = "ID"
$wsT.Range("B1").Value = "Name"
$ws1.Range("B1").Value = "Name"

# Topics data
$wsT.Range("A2").Value = 1
$wsT.Range("B2").Value = "Neurology"
$wsT.Range("A3").Value = 2
$wsT.Range("B3").Value = "Cardiology"

# Articles_Topics initial data rows (2 rows) -- Description is column D at this point
$wsAT.Range("A2").Value = 1
$wsAT.Range("B2").Value = 1
$wsAT.Range("C2").Value = 1
$wsAT.Range("D2").Value = "Neurology professor mentioned"

$wsAT.Range("A3").Value = 2
$wsAT.Range("B3").Value = 2
$wsAT.Range("C3").Value = 1
$wsAT.Range("D3").Value = "Marshall University Neurology department mentioned"

# --- Step 3: Add new article row 4 in Articles ---
$ws1.Range("A4").Value = 3
$ws1.Range("C4").Value = "http://www.mansfieldnewsjournal.com/story/news/local/2016/04/02/regions-hospitals-working-recruit-new-doctors/82511806/"
$ws1.Range("B4").Value = "Region’s hospitals working to recruit new doctors"

# --- Step 4: insert Strength column into Articles_Topics (before Description) ---
$wsAT.Columns.Item(4).Insert()
$wsAT.Range("D1").Value = "Strength"
$wsAT.Range("D2").Value = 4
$wsAT.Range("D3").Value = 4

# --- Step 5: add new rows 4, 5 to Articles_Topics ---
$wsAT.Range("A4").Value = 3
$wsAT.Range("B4").Value = 3
$wsAT.Range("C4").Value = 1
$wsAT.Range("D4").Value = 2
$wsAT.Range("E4").Value = "Neurology physicians mentioned"

$wsAT.Range("A5").Value = 4
$wsAT.Range("B5").Value = 3
$wsAT.Range("C5").Value = 2
$wsAT.Range("D5").Value = 2
$wsAT.Range("E5").Value = "Cardiology physicians mentioned"

Write-Host "done"
